# Update "Datos actualizados" timestamp in A1 (08:16 -> 08:46)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 08:46"

# Update "Muertes" (E column) counts from 11 to 16 for the affected provinces
$rows = @(29, 43, 54, 56, 58, 61, 62)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = 16
}
